$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Renumber the "id" column (A) for rows 4..70: close the gap where id=52 used
# to be skipped, so that row r ends up holding id = r-1.
for ($r = 4; $r -le 70; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}

# Add the new data row 71 (id=70, F / CE / socialmente / 1.64 / 19)
$ws.Cells.Item(71, 1).Value = 70
$ws.Cells.Item(71, 2).Value = "F"
$ws.Cells.Item(71, 3).Value = "CE"
$ws.Cells.Item(71, 4).Value = "socialmente"
$ws.Cells.Item(71, 5).Value = 1.64
$ws.Cells.Item(71, 6).Value = 19

# Move the selection to reflect where the user ended up after editing.
$ws.Range("E72").Select()
